$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B109").Value = "SETOR DE CLUBES ESPORTIVOS E ESTÁDIOS SUL"
$ws.Range("B114").Value = "SETOR DE MANSÕES ISOLADAS"
$ws.Range("B170").Value = "SETOR DE MANSÕES LAGO"
$ws.Range("B237").Value = "SETOR HABITACIONAL COLETIVO ÁREA OCTOGONAL"
$ws.Range("B280").Value = "SETOR HABITACIONAL JUSCELINO KUBITSCHEK"
